$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the two new parameter values to the "focus_email_field" row (row 5)
$ws.Range("C5").Value = "device_type"
$ws.Range("D5").Value = "user_agent"

# Move the active selection to D6, matching the saved cursor position
$ws.Range("D6").Select()
